$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceNone = 0
$vt = [char]11  # vertical tab represents a textWrapping line break in Range.Text

$rng0 = $d.Content
$found0 = $rng0.Find.Execute("ParentText is a chatbot service that helps you complete your parenting goals using a curriculum designed by Parenting for Lifelong Health with UNICEF and the World Health Organization, and tested all over the world. This programme works! ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found0) { $rng0.Text = "I-ParentText yinkonzo yekusebentisa i-chatbot lekusita nguphumelelisa migomo yakho yekuba ngumtali ngekusebentisa luhlelo lwekufundza lolwakwe yi-Parenting for Lifelong Health ne-UNICEF kanye ne-World Health Organization, futsi lwahlolwa emhlabeni wonkhe. Loluhlelo luyasebenta! " } else { Write-Host "NOT FOUND simple #0: ParentText is a chatbot service that helps you complete your parenting goals using a curriculum designed by Parenting for Lifelong Health with UNICEF and the World Health Organization, and tested all over the world. This programme works! " }

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Being here shows how much you care about providing the best support for your teen. Ngiyanihalalisela!", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found1) { $rng1.Text = "Kuba lapha kukhombisa kutsi ukhatsateke kangakanani ngekusita umntfwana wakho loseminyakeni yekutfomba. Ngiyanihalalisela!" } else { Write-Host "NOT FOUND simple #1: Being here shows how much you care about providing the best support for your teen. Ngiyanihalalisela!" }

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Remember: it is what you do with your teen that will make a difference. ParentText will provide you with tips and skills to help you with your relationship with your teen, but it is up to you to put these tips into practice!", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found2) { $rng2.Text = "Loko lokwentako nemntfwana wakho ngiko lokutakwenta umehluko. I-`"ParentText`" itakuniketa teluleko kanye nemakhono langakusita ebuhlotjeni bakho nemntfwana wakho, kodvwa kukuwe kutsi utisebentise leteluleko!" } else { Write-Host "NOT FOUND simple #2: Remember: it is what you do with your teen that will make a difference. ParentText will provide you with tips and skills to help you with your relationship with your teen, but it is up to you to put these tips into practice!" }

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("I’m Ayanda, your guide. I may look like a human, but I’m actually a robot produced by Parenting for Lifelong Health and UNICEF to help you learn. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found3) { $rng3.Text = "Ngingu-Ayanda, umsiti wakho. Ngingabukeka njengemuntfu, kodvwa ecinisweni ngiyirobothi leyentiwe yi-Parenting for Lifelong Health kanye ne-UNICEF kute ikusite ufundze. " } else { Write-Host "NOT FOUND simple #3: I’m Ayanda, your guide. I may look like a human, but I’m actually a robot produced by Parenting for Lifelong Health and UNICEF to help you learn. " }

$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Today, I’m going to explain how to use ParentText. Together we will review: ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found4) { $rng4.Text = "Lamuhla, ngitawuchaza indlela yeku sebentisa i-ParentText. Sitawuhlola. Sitawuhlola ndzawonye: " } else { Write-Host "NOT FOUND simple #4: Today, I’m going to explain how to use ParentText. Together we will review: " }

$rng5 = $d.Content
$found5 = $rng5.Find.Execute("How to earn your Positive Parenting Trophy", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found5) { $rng5.Text = "Indlela yekutfola Umklomelo Wekuba Ngumtali Lomuhle" } else { Write-Host "NOT FOUND simple #5: How to earn your Positive Parenting Trophy" }

$rng6 = $d.Content
$found6 = $rng6.Find.Execute("How to make progress in your parenting goals", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found6) { $rng6.Text = "Indlela yekutfutfukisa migomo yakho yekuba ngumtali" } else { Write-Host "NOT FOUND simple #6: How to make progress in your parenting goals" }

$rng7 = $d.Content
$found7 = $rng7.Find.Execute("How to track your progress", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found7) { $rng7.Text = "Indlela Yekulandzelela Intfutfuko yakho" } else { Write-Host "NOT FOUND simple #7: How to track your progress" }

$rng8 = $d.Content
$found8 = $rng8.Find.Execute("How to get help with this course", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found8) { $rng8.Text = "Ungalutfola njani lusito ngaleliklasi" } else { Write-Host "NOT FOUND simple #8: How to get help with this course" }

$rng9 = $d.Content
$found9 = $rng9.Find.Execute("Accessing support to troubleshoot common parenting challenges, and", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found9) { $rng9.Text = "Kutfola lusito lwekusombulula tinkinga letivamile tekuba ngumtali, kanye" } else { Write-Host "NOT FOUND simple #9: Accessing support to troubleshoot common parenting challenges, and" }

$rng10 = $d.Content
$found10 = $rng10.Find.Execute("Resources available to you in an emergency or crisis. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found10) { $rng10.Text = "Tinsita longatisebentisa uma kunesimo lesiphutfumako nobe kunenkinga. " } else { Write-Host "NOT FOUND simple #10: Resources available to you in an emergency or crisis. " }

$rng11 = $d.Content
$found11 = $rng11.Find.Execute("Earning Your Positive Parenting Trophy", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found11) { $rng11.Text = "Indlela Yekutfola Umklomelo Wekuba Ngumtali Lomuhle" } else { Write-Host "NOT FOUND simple #11: Earning Your Positive Parenting Trophy" }

$rng12 = $d.Content
$found12 = $rng12.Find.Execute("Completing a Goal", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found12) { $rng12.Text = "Kufinyelela Umgomo" } else { Write-Host "NOT FOUND simple #15: Completing a Goal" }

$rng13 = $d.Content
$found13 = $rng13.Find.Execute("Tracking Progress", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found13) { $rng13.Text = "Kulandzelela intfutfuko" } else { Write-Host "NOT FOUND simple #20: Tracking Progress" }

$rng14 = $d.Content
$found14 = $rng14.Find.Execute("Menu ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found14) { $rng14.Text = "I-Menyu " } else { Write-Host "NOT FOUND simple #25: Menu " }

$rng15 = $d.Content
$found15 = $rng15.Find.Execute("The menu contains other features that might help you, too. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found15) { $rng15.Text = "Lemenyu inetintfo letingakusita. " } else { Write-Host "NOT FOUND simple #26: The menu contains other features that might help you, too. " }

$rng16 = $d.Content
$found16 = $rng16.Find.Execute("In addition to tracking your progress, you can also: ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found16) { $rng16.Text = "Ngetulu kwekulandzelela inchubekelembili yakho, ungaphindze: " } else { Write-Host "NOT FOUND simple #27: In addition to tracking your progress, you can also: " }

$rng17 = $d.Content
$found17 = $rng17.Find.Execute("Share ParentText with a friend and help them enroll. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found17) { $rng17.Text = "Cocelenani nge-ParentText nemngani wakho futsi umsite abhalise. " } else { Write-Host "NOT FOUND simple #28: Share ParentText with a friend and help them enroll. " }

$rng18 = $d.Content
$found18 = $rng18.Find.Execute("Change your settings, like how you receive the messages, when you receive notifications, or adjust information about yourself and your teen to get the best support. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found18) { $rng18.Text = "Shintja amasethingi akho, njengendlela lowatfola ngayo umlayeto, sikhatsi lowatiswa ngaso, nobe ulungise lwati ngawe kanye newemntfwana wakho kute utfole lusito lolufanele. " } else { Write-Host "NOT FOUND simple #29: Change your settings, like how you receive the messages, when you receive notifications, or adjust information about yourself and your teen to get the best support. " }

$rng19 = $d.Content
$found19 = $rng19.Find.Execute("Access a list of activities that you can complete with your teen to build your relationship. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found19) { $rng19.Text = "Bhala luhla lwetintfo longatenta nemntfwanakho kute wakhe buhlobo benu. " } else { Write-Host "NOT FOUND simple #30: Access a list of activities that you can complete with your teen to build your relationship. " }

$rng20 = $d.Content
$found20 = $rng20.Find.Execute("Review this onboarding guide and receive support navigating ParentText. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found20) { $rng20.Text = "Hlola sicondziso sekungena bese utfola lusito lwengusebentisa iParentText. " } else { Write-Host "NOT FOUND simple #31: Review this onboarding guide and receive support navigating ParentText. " }

$rng21 = $d.Content
$found21 = $rng21.Find.Execute("And get help troubleshooting difficult challenges with your teen. Let's learn more about this feature now. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found21) { $rng21.Text = "Futsi tfole lusito lwekusombulula tinkinga letimatima umntfwana wakho lahlangabetana nato. Asesifundze kabanti ngalendzaba nyalo. " } else { Write-Host "NOT FOUND simple #32: And get help troubleshooting difficult challenges with your teen. Let's learn more about this feature now. " }

$rng22 = $d.Content
$found22 = $rng22.Find.Execute("Troubleshooting", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found22) { $rng22.Text = "Kulungisa tinkinga" } else { Write-Host "NOT FOUND simple #33: Troubleshooting" }

$rng23 = $d.Content
$found23 = $rng23.Find.Execute("Parenting can be difficult. Though challenges feel unique to you, they are often more common than you think. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found23) { $rng23.Text = "Kuba ngumtali kungaba matima. Kuba ngumtali kungaba matima. Ngisho nobe tinkinga utiva tihlukile kuwe, kodvwa tivame kakhulu kunaloko locabanga kutsi tiyenteka. " } else { Write-Host "NOT FOUND simple #34: Parenting can be difficult. Though challenges feel unique to you, they are often more common than you think. " }

$rng24 = $d.Content
$found24 = $rng24.Find.Execute("As you begin to achieve goals in this programme, I will check in on how things are going with your teen. If they aren’t going well, I might offer support. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found24) { $rng24.Text = "Njengobe ucala kufinyelela imigomo yakho kuloluhlelo, ngitawubuya ngitewubuta kutsi kuhamba njani umntfwana wakho. Ngingase ngibanikete lusito. " } else { Write-Host "NOT FOUND simple #35: As you begin to achieve goals in this programme, I will check in on how things are going with your teen. If they aren’t going well, I might offer support. " }

$rng25 = $d.Content
$found25 = $rng25.Find.Execute("When you share with me your challenges, I will offer practical solutions to help you succeed. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found25) { $rng25.Text = "Nangabe nicoca nami ngetinkinga leninato, ngitaniniketa tisombululo letitawusita niphumelele. " } else { Write-Host "NOT FOUND simple #36: When you share with me your challenges, I will offer practical solutions to help you succeed. " }

$rng26 = $d.Content
$found26 = $rng26.Find.Execute("You don’t have to wait on me to offer support, though. You can also access troubleshooting support through the Main Menu at any time. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found26) { $rng26.Text = "Kodvwa, akudzingeki ulindzele mina kutsi ngikusekele. Ungaphindze utfole lusito lwekulungisa tinkinga ngu-Main Menu nobe nini. " } else { Write-Host "NOT FOUND simple #37: You don’t have to wait on me to offer support, though. You can also access troubleshooting support through the Main Menu at any time. " }

$rng27 = $d.Content
$found27 = $rng27.Find.Execute("Emergencies and Crisis Support", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found27) { $rng27.Text = "Lusito lwetimo letiphutfumako kanye Netinkinga" } else { Write-Host "NOT FOUND simple #38: Emergencies and Crisis Support" }

$rng28 = $d.Content
$found28 = $rng28.Find.Execute("If you need information about resources in your community to address family violence, sexual violence, mental health, or other emergencies, you can message HELP at any time and receive contact details for people who may be able to help. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found28) { $rng28.Text = "Nangabe udzinga lwati ngetinsita emmangweni wakini tekubhekana nebudlova basemndenini, budlova lobentiwa ngekwelicasi, kuphatfwa kabi kwengcondvo, nobe letinye timo letiphutfumako, ungabhala umlayeto ku-LUSITO nobe kunini futsi utfole imininingwane yekutsintsana nebantfu labangakusita. " } else { Write-Host "NOT FOUND simple #39: If you need information about resources in your community to address family violence, sexual violence, mental health, or other emergencies, you can message HELP at any time and receive contact details for people who may be able to help. " }

$rng29 = $d.Content
$found29 = $rng29.Find.Execute("Your information here is safe: Nothing will be shared without your permission and will not be sold for profit. The messages you send are encrypted and locked in a secure server. However, keep in mind that if someone finds and unlocks your phone, they may be able to scroll through your messages to see what you have typed. If you send sensitive information, and you are worried, be sure to delete the messages from your phone. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found29) { $rng29.Text = "Lwati lwakho luvikelekile: Kute lutfo lolutawudluliselwa ngaphandle kwemvume yakho futsi ngeke kutsengiswe ngalo kute kuzuze. Umlayeto lowutfumelako ubhaliwe futsi ugcinwe endzaweni levikelekile. Nobe kunjalo, khumbula kutsi nangabe lotsite atfola lucingo lwakho futsi aluvule, angakhona kubuka imilayeto yakho kute abone loko bhalile. Nangabe utfumela umniningwane lobalulekile futsi ukhatsatekile, ciniseka kutsi uyawususa lomlayeto elucingweni lwakho. " } else { Write-Host "NOT FOUND simple #40: Your information here is safe: Nothing will be shared without your permission and will not be sold for profit. The messages you send are encrypted and locked in a secure server. However, keep in mind that if someone finds and unlocks your phone, they may be able to scroll through your messages to see what you have typed. If you send sensitive information, and you are worried, be sure to delete the messages from your phone. " }

$rng30 = $d.Content
$found30 = $rng30.Find.Execute("Thank you so much for listening! We hope you enjoy your ParentText journey and make the most out of it! You can access this video any time via the Main Menu. ", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", $wdReplaceNone)
if ($found30) { $rng30.Text = "Ngiyabonga kakhulu ngekungilalela! Siyetsemba kutsi utawujabulela luhambo lwakho lweParentText futsi ulisebentise kahle! Ungayi tfola le-video nobe nini ngu-main menyu. " } else { Write-Host "NOT FOUND simple #41: Thank you so much for listening! We hope you enjoy your ParentText journey and make the most out of it! You can access this video any time via the Main Menu. " }

$gp27 = $d.Paragraphs.Item(27)
$gp27.Range.Text = "Kwekucala, indlela yekucedzela sifundvo kute utfole umklomelo wakho wekuba ngumtali lokahle. " + $vt + $vt + " Ngekuya ngesimo sakho, lesifundvo sihlukaniswe saba imigomo lengu-8 nobe lengu-9 imigomo yekukhulisa lengakusita utfutfukise buhlobo bakho nemntfwana wakho futsi umsite aphumelele. Umgomo ngamunye utsatsa emalanga lamabili kuya kulamane kutsi uwucedzele, futsi lilanga ngalinye liletsa sifundvo lesisha. " + $vt + $vt + " Ungakhetsa kutsi ngumiphi imigomo lofuna kuyenta kucala, kodvwa umgomo ngamunye kufanele ucedvwe kute ucedzele sifundvo futsi utfole umklomelo wakho wekuba ngumtali lomuhle."

$gp31 = $d.Paragraphs.Item(31)
$gp31.Range.Text = " Nyalo-ke, asesifundze kutsi singaticedza njani tinhloso taloluhlelo. " + $vt + $vt + "Masinyane, utawucelwa kutsi ukhetse umgomo wakho wekucala wekuba ngumtali. Nasewukhetse umgomo, utawucedzela tifundvo kute utfole emakhono lamasha. " + $vt + $vt + " Onkhe emalanga kunemakhono lamasha. Emakhono lamanyenti atsatsa emaminitsi langaphasi kwalangu-5. Nangabe ungakhoni kucedza likhono lolibelwe, ngitawubuta kutsi uyafuna yini kucedza lelikhono ngelilanga lelilandzelako. " + $vt + $vt + " Kumele cedze onkhe emakhono kute utfole ibheji. Nasewuyifezile migomo yakho yekuba ngumtali, ungakhetsa lenye. Uma sewutfole onkhe emabheji akho, utawube sewucedzile kufundza futsi utawutfola umklomelo wekukhulisa bantfwana ngendlela lekahle."

$gp35 = $d.Paragraphs.Item(35)
$gp35.Range.Text = "Njengobe uchubeka nesifundvo ngasinye, utawutfola lwati lolubonakala ngalendlela. Letibonakaliso tikutjela kutsi sewuhambile kangakanani esifundvweni. " + $vt + $vt + " Nangabe ufuna kubona kutsi sewutfutfuke kanganani ekufinyeleleni imigomo yakho yekuba ngumtali, ungahlola inchubekelembili loyitfolile ngekusebentisa i-Main Menyu. Kute ufinyelele imenyu, bhala `"Imenyu`" nobe nini. " + $vt + $vt + " Inketfo yekucala ngu-menyu ibhalwe `"Landzelela inchubekelembili yami`". " + $vt + $vt + " Lapha ungabona inchubekelembili yakho, uphindze ubuke imigomo loyifinyelele kanye naleyo lengakacedvwa."

